$d = $word.ActiveDocument

function SplitAt($pos) {
    $r = $d.Range($pos, $pos)
    $bmName = "TmpSplit_" + $pos
    $d.Bookmarks.Add($bmName, $r)
    $d.Bookmarks($bmName).Delete()
}

# Locate the start of the target paragraph ("I often hear disputes ...")
$findRng = $d.Content
$findRng.Find.Execute("I often hear disputes")
$paraStart = $findRng.Start

# --- Step 1: perform the 3 text insertions, from rightmost to leftmost so
# earlier offsets stay valid. ---

# Insertion C: " software. I think that" -> " software. I think  that" (extra space before "that")
$insC = $paraStart + 265
$rC = $d.Range($insC, $insC)
$rC.InsertAfter(" ")

# Insertion B: " in 2011, it turned out that open source " -> insert " the" after "that"
$insB = $paraStart + 180
$rB = $d.Range($insB, $insB)
$rB.InsertAfter(" the")

# Insertion A: "According to the results of researches..." -> insert " the" after "of"
$insA = $paraStart + 97
$rA = $d.Range($insA, $insA)
$rA.InsertAfter(" the")

"Insertions done"

# --- Step 2: re-establish every run boundary required by the final layout. ---
# (Positions are absolute, computed on the POST-insertion text, i.e. already
# account for the 4 + 4 + 1 = 9 extra characters inserted above.)
$boundaries = @(1189, 1206, 1233, 1237, 1268, 1276, 1284, 1292, 1320, 1324, 1337, 1345, 1391, 1404, 1410, 1411, 1415, 1501, 1505)
foreach ($b in $boundaries) {
    SplitAt($b)
}

"Splits done"
